$wb = $excel.ActiveWorkbook

# Rename sheets: swap "Users_input" and "Users_input (backup)" (note the
# target name gains a space: "Users_input (back up)").
$wsBackup = $wb.Worksheets.Item("Users_input")
$wsBackup.Name = "Users_input (back up)"

$wsMain = $wb.Worksheets.Item("Users_input (backup)")
$wsMain.Name = "Users_input"

# Make "Users_input" (the renamed former backup sheet) the active tab and
# update its selection.
$wsMain.Activate()
[void]$wsMain.Range("C26").Select()
